# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 502 (pushing existing rows 502:540 down to 503:541).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting a new blank row at 502.
$ws.Rows.Item(502).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A502").Value = 7
$ws.Range("B502").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C502").Value = 'Ñuble'
$ws.Range("D502").Value = 45223
$ws.Range("E502").Value = 16
$ws.Range("F502").Value = 100112003
$ws.Range("G502").Value = 'Ajo'
$ws.Range("H502").Value = 'Chino'
$ws.Range("I502").Value = 'Primera'
$ws.Range("J502").Value = 50
$ws.Range("K502").Value = 22000
$ws.Range("L502").Value = 22000
$ws.Range("M502").Value = 22000
$ws.Range("N502").Value = '$/malla 10 kilos'
$ws.Range("O502").Value = 'China'
$ws.Range("P502").Value = 2200
$ws.Range("Q502").Value = 10
$ws.Range("R502").Value = 'Hortaliza'
